# Add the new "V2.01 / Level LED와 PWR LED 동기 시킴" revision row to the
# "Revision" history sheet (F/W change: LEVEL LED synced with PWR LED on power off).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revision")

# Row 9 was a blank template row below the "V2.00_1.0sec" entry (row 8) -
# fill it in with the new revision's date, version and description.
$ws.Range("B9").Value = 43565
$ws.Range("C9").Value = "V2.01"
$ws.Range("D9").Value = "Level LED와 PWR LED 동기 시킴"

# Leave the selection where the author left it after typing the new row.
$ws.Range("D10").Select()
